# Replicates the authored changes to Varianty.xlsx:
#  - sheet renamed from "Sheet1" to "List1"
#  - active selection on the sheet moved from D9 to A13
#
# (The window-geometry / revision-session attributes seen in the raw XML
# diff are Excel-desktop session artifacts recorded at save time and are
# not part of the workbook's editable content/object model.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the only worksheet.
$ws.Name = "List1"

# Move the selection/active cell to A13.
$ws.Range("A13").Select()
